# Generate Report for Archive
# 1. Status text "Ready for handoff" -> "In Translation" on every sheet that
#    tracks handoff/translation status (Overview's zh-cn/de-de status columns,
#    plus the per-locale "Status" column on the zh-cn and de-de sheets).
# 2. The Status column narrows (AutoFit-style) to fit the shorter text, which
#    also drives the matching zh-cn/de-de columns on the Overview sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# Narrow the Status-related columns to match the shorter "In Translation" text.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
